$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 32260356
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 35716610
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 35716610
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -35716960
$ws.Range("H137").Value = 1117.1818
$ws.Range("I137").Value = 922.44446
$ws.Range("J137").Value = 1993.5
$ws.Range("K137").Value = 2767.33338
$ws.Range("L137").Value = 5980.5
$ws.Range("M137").Value = -217.33338
$ws.Range("N137").Value = -11080.5
$ws.Range("H141").Value = 2742.7778
$ws.Range("I141").Value = 2266.5625
$ws.Range("J141").Value = 6552.5
$ws.Range("K141").Value = 6799.6875
$ws.Range("L141").Value = 19657.5
$ws.Range("M141").Value = -1619.6875
$ws.Range("N141").Value = -30017.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2020828.5
$ws.Range("I32").Value = 3799.2036
$ws.Range("K32").Value = 3799.2036
$ws.Range("M32").Value = -3512.2036
$ws.Range("H70").Value = 34000
$ws.Range("J70").Value = 34000
$ws.Range("L70").Value = 34000
$ws.Range("N70").Value = -34540
$ws.Range("H73").Value = 34000
$ws.Range("J73").Value = 34000
$ws.Range("L73").Value = 34000
$ws.Range("N73").Value = -35872
$ws.Range("H74").Value = 831.7222
$ws.Range("I74").Value = 783.08887
$ws.Range("J74").Value = 1074.8889
$ws.Range("K74").Value = 783.08887
$ws.Range("L74").Value = 1074.8889
$ws.Range("M74").Value = 90.91112999999996
$ws.Range("N74").Value = -2822.8889
$ws.Range("H77").Value = 831.7222
$ws.Range("I77").Value = 783.08887
$ws.Range("J77").Value = 1074.8889
$ws.Range("K77").Value = 3915.44435
$ws.Range("L77").Value = 5374.4445
$ws.Range("M77").Value = 452.5556499999998
$ws.Range("N77").Value = -14110.4445
$ws.Range("H132").Value = 1701.2759
$ws.Range("I132").Value = 1170.4706
$ws.Range("J132").Value = 2453.25
$ws.Range("K132").Value = 3511.4118
$ws.Range("L132").Value = 7359.75
$ws.Range("M132").Value = -981.4118000000003
$ws.Range("N132").Value = -12419.75
$ws.Range("H139").Value = 45710
$ws.Range("J139").Value = 45710
$ws.Range("L139").Value = 45710
$ws.Range("N139").Value = -55990

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 499.66666
$ws.Range("I22").Value = 499.66666
$ws.Range("K22").Value = 499.66666
$ws.Range("M22").Value = -326.66666
$ws.Range("H132").Value = 219666.67
$ws.Range("J132").Value = 219666.67
$ws.Range("L132").Value = 219666.67
$ws.Range("N132").Value = -229786.67
$ws.Range("H134").Value = 29174.744
$ws.Range("I134").Value = 3764.3635
$ws.Range("J134").Value = 62058.766
$ws.Range("K134").Value = 11293.0905
$ws.Range("L134").Value = 186176.298
$ws.Range("M134").Value = -8758.0905
$ws.Range("N134").Value = -191246.298
$ws.Range("H135").Value = 60639.5
$ws.Range("J135").Value = 60639.5
$ws.Range("L135").Value = 60639.5
$ws.Range("N135").Value = -70779.5
$ws.Range("H137").Value = 55338.46
$ws.Range("J137").Value = 55338.46
$ws.Range("L137").Value = 55338.46
$ws.Range("N137").Value = -65538.45999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2834.0476
$ws.Range("I31").Value = 2940.75
$ws.Range("J31").Value = 700
$ws.Range("K31").Value = 2940.75
$ws.Range("L31").Value = 700
$ws.Range("M31").Value = -2645.75
$ws.Range("N31").Value = -1290
$ws.Range("H34").Value = 2834.0476
$ws.Range("I34").Value = 2940.75
$ws.Range("J34").Value = 700
$ws.Range("K34").Value = 2940.75
$ws.Range("L34").Value = 700
$ws.Range("M34").Value = -2738.75
$ws.Range("N34").Value = -1104
$ws.Range("H86").Value = 7441.6665
$ws.Range("I86").Value = 3460
$ws.Range("J86").Value = 10285.714
$ws.Range("K86").Value = 3460
$ws.Range("L86").Value = 10285.714
$ws.Range("M86").Value = -2337
$ws.Range("N86").Value = -12531.714
$ws.Range("H89").Value = 7441.6665
$ws.Range("I89").Value = 3460
$ws.Range("J89").Value = 10285.714
$ws.Range("K89").Value = 17300
$ws.Range("L89").Value = 51428.57
$ws.Range("M89").Value = -11684
$ws.Range("N89").Value = -62660.57
$ws.Range("H134").Value = 4685.615
$ws.Range("I134").Value = 4164.125
$ws.Range("J134").Value = 5520
$ws.Range("K134").Value = 12492.375
$ws.Range("L134").Value = 16560
$ws.Range("M134").Value = -9957.375
$ws.Range("N134").Value = -21630

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 551.8
$ws.Range("I121").Value = 364.75
$ws.Range("J121").Value = 1300
$ws.Range("K121").Value = 1094.25
$ws.Range("L121").Value = 3900
$ws.Range("M121").Value = 215.75
$ws.Range("N121").Value = -6520
$ws.Range("H122").Value = 345671.75
$ws.Range("J122").Value = 770439
$ws.Range("L122").Value = 6933951
$ws.Range("N122").Value = -6938851
$ws.Range("H123").Value = 10000
$ws.Range("J123").Value = 10000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -34900
$ws.Range("H131").Value = 850.64
$ws.Range("I131").Value = 515.36365
$ws.Range("J131").Value = 892.0787
$ws.Range("K131").Value = 1546.09095
$ws.Range("L131").Value = 2676.2361
$ws.Range("M131").Value = 3493.90905
$ws.Range("N131").Value = -12756.2361

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1257.3704
$ws.Range("I22").Value = 831.4286
$ws.Range("J22").Value = 1406.45
$ws.Range("K22").Value = 831.4286
$ws.Range("L22").Value = 1406.45
$ws.Range("M22").Value = -536.4286
$ws.Range("N22").Value = -1996.45
$ws.Range("H27").Value = 1257.3704
$ws.Range("I27").Value = 831.4286
$ws.Range("J27").Value = 1406.45
$ws.Range("K27").Value = 831.4286
$ws.Range("L27").Value = 1406.45
$ws.Range("M27").Value = -724.4286
$ws.Range("N27").Value = -1620.45

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6100
$ws.Range("I62").Value = 4614.2856
$ws.Range("J62").Value = 9566.666999999999
$ws.Range("K62").Value = 4614.2856
$ws.Range("L62").Value = 9566.666999999999
$ws.Range("M62").Value = -3990.2856
$ws.Range("N62").Value = -10814.667
$ws.Range("H65").Value = 6100
$ws.Range("I65").Value = 4614.2856
$ws.Range("J65").Value = 9566.666999999999
$ws.Range("K65").Value = 23071.428
$ws.Range("L65").Value = 47833.335
$ws.Range("M65").Value = -19951.428
$ws.Range("N65").Value = -54073.335
